$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.137.29"
$ws.Range("E2").Value = "  +0.58%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.770.71"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB (numeric-looking price, force text so it isn't coerced to a number)
$ws.Range("D5").Value = "'593.32"
$ws.Range("E5").Value = "  +0.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'167.28"

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.768.87"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.13%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.31%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.40"
$ws.Range("E11").Value = "  -1.33%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -0.56%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -1.51%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'36.05"
$ws.Range("E14").Value = "  -0.81%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.402.32"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.752.04"
$ws.Range("E16").Value = "  +0.96%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "68.125.36"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "'17.83"
$ws.Range("E18").Value = "  -2.22%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.43%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "'6.97"
$ws.Range("E20").Value = "  -0.85%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'10.78"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'463.64"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -0.66%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +9.58%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "'83.81"
$ws.Range("E25").Value = "  +1.35%  "

# Row 26 - Fetch.AI
$ws.Range("D26").Value = "'2.17"
$ws.Range("E26").Value = "  -1.63%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'11.81"
$ws.Range("E27").Value = "  -1.64%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -0.91%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.13%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  -0.51%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'29.85"
$ws.Range("E32").Value = "  +0.76%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "'2.15"
$ws.Range("E33").Value = "  -3.66%  "

# Row 34 & 35 - Aptos / Binance-PegBSC-USD swap ranking positions
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'9.11"
$ws.Range("E35").Value = "  +0.80%  "

# Row 36 - RenzoRestakedETH
$ws.Range("D36").Value = "3.723.71"
$ws.Range("E36").Value = "  +1.23%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -1.02%  "

# Row 38 - dogwifhat
$ws.Range("D38").Value = "'3.43"
$ws.Range("E38").Value = "  +0.07%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +0.27%  "

# Row 40 - Mantle
$ws.Range("E40").Value = "  +1.29%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  -0.04%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.05%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - Arweave
$ws.Range("D44").Value = "'44.25"
$ws.Range("E44").Value = "  +16.20%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -2.18%  "

# Row 46 - OKB
$ws.Range("D46").Value = "'46.92"
$ws.Range("E46").Value = "  +3.34%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -0.83%  "

# Row 48 - Cosmos
$ws.Range("E48").Value = "  -1.71%  "

# Row 49 - Monero
$ws.Range("D49").Value = "'145.36"
$ws.Range("E49").Value = "  +1.41%  "

# Row 50 - Bittensor
$ws.Range("D50").Value = "'387.83"
$ws.Range("E50").Value = "  -0.63%  "

# Row 51 - Maker
$ws.Range("D51").Value = "2.781.72"
$ws.Range("E51").Value = "  +4.05%  "
